# Reverted default group numbers and metadata so they are consistent with
# the other default paths and data.
#
# Updates the "UnitMass" (column C) values for the +loading block (rows 2-21)
# and the -loading block (rows 23-42) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "C2" = 23;   "C3" = 125;  "C4" = 83;   "C5" = 69;   "C6" = 68;
    "C7" = 82;   "C8" = 46;   "C9" = 32;   "C10" = 28;  "C11" = 126;
    "C12" = 81;  "C13" = 95;  "C14" = 15;  "C15" = 45;  "C16" = 42;
    "C17" = 61;  "C18" = 30;  "C19" = 97;  "C20" = 27;  "C21" = 54;
    "C23" = 39;  "C24" = 71;  "C25" = 55;  "C26" = 41;  "C27" = 43;
    "C28" = 102; "C29" = 40;  "C30" = 175; "C31" = 57;  "C32" = 149;
    "C33" = 59;  "C34" = 115; "C35" = 103; "C36" = 231; "C37" = 74;
    "C38" = 112; "C39" = 77;  "C40" = 53;  "C41" = 65;  "C42" = 73;
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
